$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 110.333336
$ws.Range("I4").Value = 110.333336
$ws.Range("K4").Value = 110.333336
$ws.Range("M4").Value = 3.666663999999997
$ws.Range("H76").Value = 6175795
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 7939393.5
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 7939393.5
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -7940023.5
$ws.Range("H79").Value = 6175795
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 7939393.5
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 7939393.5
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -7941577.5
$ws.Range("H86").Value = 5831.636
$ws.Range("I86").Value = 1308
$ws.Range("J86").Value = 11260
$ws.Range("K86").Value = 1308
$ws.Range("L86").Value = 11260
$ws.Range("M86").Value = -185
$ws.Range("N86").Value = -13506
$ws.Range("H89").Value = 5831.636
$ws.Range("I89").Value = 1308
$ws.Range("J89").Value = 11260
$ws.Range("K89").Value = 6540
$ws.Range("L89").Value = 56300
$ws.Range("M89").Value = -924
$ws.Range("N89").Value = -67532
$ws.Range("H96").Value = 1586.1538
$ws.Range("I96").Value = 1187.75
$ws.Range("J96").Value = 1763.2222
$ws.Range("K96").Value = 3563.25
$ws.Range("L96").Value = 5289.6666
$ws.Range("M96").Value = -2190.25
$ws.Range("N96").Value = -8035.6666
$ws.Range("H100").Value = 1675.25
$ws.Range("I100").Value = 1067.3334
$ws.Range("K100").Value = 1067.3334
$ws.Range("M100").Value = -526.3334
$ws.Range("H129").Value = 763.4107
$ws.Range("J129").Value = 797.38464
$ws.Range("L129").Value = 2392.15392
$ws.Range("N129").Value = -12392.15392
$ws.Range("H132").Value = 3350.6453
$ws.Range("I132").Value = 3335.963
$ws.Range("K132").Value = 10007.889
$ws.Range("M132").Value = -7477.889000000001
$ws.Range("H137").Value = 79897.92
$ws.Range("I137").Value = 96404.48
$ws.Range("J137").Value = 2867.3333
$ws.Range("K137").Value = 289213.44
$ws.Range("L137").Value = 8601.999899999999
$ws.Range("M137").Value = -286663.44
$ws.Range("N137").Value = -13701.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1882.4474
$ws.Range("I61").Value = 1549.5161
$ws.Range("K61").Value = 1549.5161
$ws.Range("M61").Value = -1337.5161
$ws.Range("H122").Value = 2000.0714
$ws.Range("I122").Value = 1813.4166
$ws.Range("J122").Value = 3120
$ws.Range("K122").Value = 5440.2498
$ws.Range("L122").Value = 9360
$ws.Range("M122").Value = -2990.2498
$ws.Range("N122").Value = -14260
$ws.Range("H136").Value = 1882.4474
$ws.Range("I136").Value = 1549.5161
$ws.Range("K136").Value = 4648.5483
$ws.Range("M136").Value = -2098.5483

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 28689.5
$ws.Range("J52").Value = 28689.5
$ws.Range("L52").Value = 28689.5
$ws.Range("N52").Value = -29215.5
$ws.Range("H80").Value = 851.3333
$ws.Range("I80").Value = 775.5454999999999
$ws.Range("J80").Value = 895.2105
$ws.Range("K80").Value = 775.5454999999999
$ws.Range("L80").Value = 895.2105
$ws.Range("M80").Value = 222.4545000000001
$ws.Range("N80").Value = -2891.2105
$ws.Range("H83").Value = 851.3333
$ws.Range("I83").Value = 775.5454999999999
$ws.Range("J83").Value = 895.2105
$ws.Range("K83").Value = 3877.7275
$ws.Range("L83").Value = 4476.0525
$ws.Range("M83").Value = 1114.2725
$ws.Range("N83").Value = -14460.0525
$ws.Range("H99").Value = 1432
$ws.Range("I99").Value = 1307.5
$ws.Range("J99").Value = 1587.625
$ws.Range("K99").Value = 1307.5
$ws.Range("L99").Value = 1587.625
$ws.Range("M99").Value = 190.5
$ws.Range("N99").Value = -4583.625
$ws.Range("H121").Value = 28689.5
$ws.Range("J121").Value = 28689.5
$ws.Range("L121").Value = 28689.5
$ws.Range("N121").Value = -32183.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1121.5454
$ws.Range("I16").Value = 737.6
$ws.Range("J16").Value = 1441.5
$ws.Range("K16").Value = 737.6
$ws.Range("L16").Value = 1441.5
$ws.Range("M16").Value = -450.6
$ws.Range("N16").Value = -2015.5
$ws.Range("H105").Value = 799.25
$ws.Range("I105").Value = 799.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 799.25
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 947.75
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 1121.5454
$ws.Range("I113").Value = 737.6
$ws.Range("J113").Value = 1441.5
$ws.Range("K113").Value = 737.6
$ws.Range("L113").Value = 1441.5
$ws.Range("M113").Value = 1432.4
$ws.Range("N113").Value = -5781.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 750.29
$ws.Range("J131").Value = 763.9382000000001
$ws.Range("L131").Value = 2291.8146
$ws.Range("N131").Value = -12371.8146

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9525.833000000001
$ws.Range("I113").Value = 12476.375
$ws.Range("J113").Value = 3624.75
$ws.Range("K113").Value = 12476.375
$ws.Range("L113").Value = 3624.75
$ws.Range("M113").Value = -10306.375
$ws.Range("N113").Value = -7964.75
$ws.Range("H122").Value = 4280.65
$ws.Range("I122").Value = 4261.4
$ws.Range("J122").Value = 4299.9
$ws.Range("K122").Value = 12784.2
$ws.Range("L122").Value = 12899.7
$ws.Range("M122").Value = -10334.2
$ws.Range("N122").Value = -17799.7

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5526.5293
$ws.Range("J22").Value = 4333
$ws.Range("L22").Value = 4333
$ws.Range("N22").Value = -4923
$ws.Range("H27").Value = 5526.5293
$ws.Range("J27").Value = 4333
$ws.Range("L27").Value = 4333
$ws.Range("N27").Value = -4547
$ws.Range("H40").Value = 3521.48
$ws.Range("I40").Value = 3332.35
$ws.Range("J40").Value = 4278
$ws.Range("K40").Value = 3332.35
$ws.Range("L40").Value = 4278
$ws.Range("M40").Value = -3196.35
$ws.Range("N40").Value = -4550
$ws.Range("H132").Value = 505101.6
$ws.Range("I132").Value = 636128.4399999999
$ws.Range("J132").Value = 7199.6
$ws.Range("K132").Value = 1908385.32
$ws.Range("L132").Value = 21598.8
$ws.Range("M132").Value = -1905855.32
$ws.Range("N132").Value = -26658.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 293.84616
$ws.Range("I100").Value = 294.54544
$ws.Range("J100").Value = 290
$ws.Range("K100").Value = 589.09088
$ws.Range("L100").Value = 580
$ws.Range("M100").Value = -48.09087999999997
$ws.Range("N100").Value = -1662
$ws.Range("H122").Value = 1199.6428
$ws.Range("I122").Value = 941.6667
$ws.Range("K122").Value = 2825.0001
$ws.Range("M122").Value = -375.0001000000002
$ws.Range("H132").Value = 1739.4584
$ws.Range("I132").Value = 1274.25
$ws.Range("K132").Value = 3822.75
$ws.Range("M132").Value = -1292.75
